$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "name" column (column A) entirely - "sampleid" becomes the
# sole identifying field, per "switching to name as only required field".
$ws.Range("A1").EntireColumn.Delete()

# Update the active selection to match the post-edit cursor position.
$ws.Range("B12").Select()
